$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to remain plain text so Excel does not
# auto-convert numeric-looking strings (e.g. "0.9988") into floating point
# numbers, which would change both the stored type and the rendered value.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.971.36'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '1.899.63'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '0.7898'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").Value = '244.53'
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("D7").Value = '0.9991'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '0.3162'
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").Value = '25.83'
$ws.Range("E9").Value = '  +1.60%  '
$ws.Range("D10").Value = '0.07333'
$ws.Range("E10").Value = '  +4.65%  '
$ws.Range("D11").Value = '0.08135'
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("D12").Value = '0.7773'
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("D13").Value = '5.529'
$ws.Range("E13").Value = '  +4.43%  '
$ws.Range("D14").Value = '94.37'
$ws.Range("E14").Value = '  +2.37%  '
$ws.Range("D15").Value = '1.813.59'
$ws.Range("E15").Value = '  -4.32%  '
$ws.Range("D16").Value = '6.267'
$ws.Range("E16").Value = '  +5.89%  '
$ws.Range("D17").Value = '29.784.64'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '14.02'
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").Value = '247.51'
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").Value = '0.000007853'
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("D21").Value = '8.179'
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").Value = '0.9987'
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").Value = '2.092.45'
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").Value = '0.9986'
$ws.Range("D25").Value = '0.1604'
$ws.Range("E25").Value = '  -2.23%  '
$ws.Range("D26").Value = '9.509'
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").Value = '163.83'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '18.86'
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").Value = '2.048'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '1.449'
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").Value = '4.505'
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("D33").Value = '0.05632'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '4.107'
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").Value = '1.254'
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").Value = '0.7571'
$ws.Range("E36").Value = '  +3.11%  '
$ws.Range("D38").Value = '2.661'
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").Value = '0.01941'
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").Value = '2.797'
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").Value = '1.149.40'
$ws.Range("E41").Value = '  +12.12%  '
$ws.Range("D42").Value = '0.4482'
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("D44").Value = '5.987'
$ws.Range("E44").Value = '  +2.85%  '
$ws.Range("D45").Value = '0.8595'
$ws.Range("E45").Value = '  +2.38%  '
$ws.Range("B46").Value = 'SynthetixNetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D46").Value = '3.177'
$ws.Range("E46").Value = '  +10.00%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.906'
$ws.Range("E47").Value = '  +3.00%  '
$ws.Range("D48").Value = '0.9991'
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").Value = '102.33'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("D51").Value = '9.795'
$ws.Range("E51").Value = '  -1.50%  '

# Restore the default "Normal" style on the price cells so no stray
# number-format style index is left behind on the saved workbook.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
